# Generate Report for Handback
# - Update the "Ready for handoff" status text (used by both zh-cn and de-de sheets'
#   row for 712274be-... .md) to "Handback transform failed".
# - Populate the "Error Detail" cell for the 712274be-... row on both language sheets
#   with a handback/handoff filename-mismatch message.
# - Widen the "Error Detail" column (column P, the 16th column) on both language
#   sheets so the new message is readable.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_dede = $wb.Worksheets.Item("de-de")

# Status text shared by the "712274be-d5bc-4510-9f48-a0d6668212c5.md" row everywhere it is
# surfaced: Overview (zh-cn/de-de columns) and the per-language sheets' Status column.
$ws_overview.Range("E3").Value = "Handback transform failed"
$ws_overview.Range("F3").Value = "Handback transform failed"
$ws_zhcn.Range("C3").Value = "Handback transform failed"
$ws_dede.Range("C3").Value = "Handback transform failed"

# New Error Detail messages for row 3 on each language sheet
$ws_zhcn.Range("P3").Value = "Handback file name: zv4x4equ.lp1 is different with handoff file name: 712274be-d5bc-4510-9f48-a0d6668212c5.22987532d7d88b3a630fb28e33d3345ff223fab5.zh-cn."
$ws_dede.Range("P3").Value = "Handback file name: zv4x4equ.lp1 is different with handoff file name: 712274be-d5bc-4510-9f48-a0d6668212c5.22987532d7d88b3a630fb28e33d3345ff223fab5.de-de."

# Widen column P (Error Detail) to fit the new text on both sheets
$ws_zhcn.Columns.Item(16).ColumnWidth = 39.15
$ws_dede.Columns.Item(16).ColumnWidth = 39.15
